$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 370.94116
$ws.Range("J58").Value = 587.5
$ws.Range("L58").Value = 1762.5
$ws.Range("N58").Value = -2062.5
$ws.Range("H64").Value = 5125.875
$ws.Range("I64").Value = 5334.3335
$ws.Range("J64").Value = 5000.8
$ws.Range("K64").Value = 5334.3335
$ws.Range("L64").Value = 5000.8
$ws.Range("M64").Value = -5086.3335
$ws.Range("N64").Value = -5496.8
$ws.Range("H67").Value = 5125.875
$ws.Range("I67").Value = 5334.3335
$ws.Range("J67").Value = 5000.8
$ws.Range("K67").Value = 5334.3335
$ws.Range("L67").Value = 5000.8
$ws.Range("M67").Value = -4476.3335
$ws.Range("N67").Value = -6716.8
$ws.Range("H132").Value = 3002.9138
$ws.Range("I132").Value = 2260.3
$ws.Range("J132").Value = 7644.25
$ws.Range("K132").Value = 6780.900000000001
$ws.Range("L132").Value = 22932.75
$ws.Range("M132").Value = -4250.900000000001
$ws.Range("N132").Value = -27992.75
$ws.Range("H137").Value = 4195.9565
$ws.Range("I137").Value = 1283.6
$ws.Range("J137").Value = 6436.231
$ws.Range("K137").Value = 3850.8
$ws.Range("L137").Value = 19308.693
$ws.Range("M137").Value = -1300.8
$ws.Range("N137").Value = -24408.693
$ws.Range("H138").Value = 3667.681
$ws.Range("I138").Value = 3595.7778
$ws.Range("J138").Value = 3684.7104
$ws.Range("K138").Value = 10787.3334
$ws.Range("L138").Value = 11054.1312
$ws.Range("M138").Value = -5647.3334
$ws.Range("N138").Value = -21334.1312

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7374.897
$ws.Range("I32").Value = 4336.1206
$ws.Range("K32").Value = 4336.1206
$ws.Range("M32").Value = -4049.1206
$ws.Range("H45").Value = 2471.5417
$ws.Range("I45").Value = 1095.9
$ws.Range("J45").Value = 3454.1428
$ws.Range("K45").Value = 1095.9
$ws.Range("L45").Value = 3454.1428
$ws.Range("M45").Value = -718.9000000000001
$ws.Range("N45").Value = -4208.1428
$ws.Range("H63").Value = 2999
$ws.Range("I63").Value = 2999
$ws.Range("K63").Value = 2999
$ws.Range("M63").Value = -2313
$ws.Range("H66").Value = 2999
$ws.Range("I66").Value = 2999
$ws.Range("K66").Value = 14995
$ws.Range("M66").Value = -11563
$ws.Range("H74").Value = 145729.58
$ws.Range("I74").Value = 168684.5
$ws.Range("K74").Value = 168684.5
$ws.Range("M74").Value = -167810.5
$ws.Range("H77").Value = 145729.58
$ws.Range("I77").Value = 168684.5
$ws.Range("K77").Value = 843422.5
$ws.Range("M77").Value = -839054.5
$ws.Range("H122").Value = 3010.7144
$ws.Range("I122").Value = 2857.6924
$ws.Range("K122").Value = 8573.0772
$ws.Range("M122").Value = -6123.0772
$ws.Range("H132").Value = 39517.742
$ws.Range("I132").Value = 47835.727
$ws.Range("J132").Value = 2918.6
$ws.Range("K132").Value = 143507.181
$ws.Range("L132").Value = 8755.799999999999
$ws.Range("M132").Value = -140977.181
$ws.Range("N132").Value = -13815.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1793.8889
$ws.Range("I5").Value = 1793.8889
$ws.Range("K5").Value = 1793.8889
$ws.Range("M5").Value = -1680.8889
$ws.Range("H86").Value = 3796
$ws.Range("I86").Value = 2500
$ws.Range("K86").Value = 2500
$ws.Range("M86").Value = -1377
$ws.Range("H89").Value = 3796
$ws.Range("I89").Value = 2500
$ws.Range("K89").Value = 12500
$ws.Range("M89").Value = -6884
$ws.Range("H134").Value = 2149.2354
$ws.Range("I134").Value = 1948.0652
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 5844.1956
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -3309.1956
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 7999.5
$ws.Range("J14").Value = 7999.5
$ws.Range("L14").Value = 7999.5
$ws.Range("N14").Value = -8339.5
$ws.Range("H31").Value = 5726.8125
$ws.Range("I31").Value = 4466.25
$ws.Range("K31").Value = 4466.25
$ws.Range("M31").Value = -4171.25
$ws.Range("H34").Value = 5726.8125
$ws.Range("I34").Value = 4466.25
$ws.Range("K34").Value = 4466.25
$ws.Range("M34").Value = -4264.25
$ws.Range("H122").Value = 1012.2727
$ws.Range("I122").Value = 1127.4445
$ws.Range("K122").Value = 3382.3335
$ws.Range("M122").Value = -932.3335000000002
$ws.Range("H132").Value = 4269.4443
$ws.Range("I132").Value = 4698.615
$ws.Range("J132").Value = 3153.6
$ws.Range("K132").Value = 14095.845
$ws.Range("L132").Value = 9460.799999999999
$ws.Range("M132").Value = -11565.845
$ws.Range("N132").Value = -14520.8
$ws.Range("H134").Value = 51012.477
$ws.Range("I134").Value = 56095.42
$ws.Range("K134").Value = 168286.26
$ws.Range("M134").Value = -165751.26

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4696.5
$ws.Range("J39").Value = 5253.143
$ws.Range("L39").Value = 15759.429
$ws.Range("N39").Value = -16347.429
$ws.Range("H55").Value = 2825
$ws.Range("J55").Value = 3175
$ws.Range("L55").Value = 9525
$ws.Range("N55").Value = -9879
$ws.Range("H132").Value = 974.1539
$ws.Range("I132").Value = 896.125
$ws.Range("J132").Value = 1099
$ws.Range("K132").Value = 8065.125
$ws.Range("L132").Value = 9891
$ws.Range("M132").Value = -5535.125
$ws.Range("N132").Value = -14951

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 20949
$ws.Range("I41").Value = 1899
$ws.Range("J41").Value = 39999
$ws.Range("K41").Value = 1899
$ws.Range("L41").Value = 39999
$ws.Range("M41").Value = -1544
$ws.Range("N41").Value = -40709
$ws.Range("H80").Value = 4373.5
$ws.Range("I80").Value = 3750.5
$ws.Range("K80").Value = 3750.5
$ws.Range("M80").Value = -2752.5
$ws.Range("H83").Value = 4373.5
$ws.Range("I83").Value = 3750.5
$ws.Range("K83").Value = 18752.5
$ws.Range("M83").Value = -13760.5
$ws.Range("H99").Value = 14893.333
$ws.Range("I99").Value = 8005.125
$ws.Range("K99").Value = 8005.125
$ws.Range("M99").Value = -5759.125
$ws.Range("H122").Value = 2239.5
$ws.Range("I122").Value = 1798.4
$ws.Range("K122").Value = 5395.200000000001
$ws.Range("M122").Value = -2945.200000000001
$ws.Range("H132").Value = 30672.945
$ws.Range("I132").Value = 42220.76
$ws.Range("K132").Value = 126662.28
$ws.Range("M132").Value = -124132.28

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8972.406999999999
$ws.Range("I7").Value = 10415.75
$ws.Range("J7").Value = 4848.5713
$ws.Range("K7").Value = 10415.75
$ws.Range("L7").Value = 4848.5713
$ws.Range("M7").Value = -10303.75
$ws.Range("N7").Value = -5072.5713
$ws.Range("H126").Value = 8972.406999999999
$ws.Range("I126").Value = 10415.75
$ws.Range("J126").Value = 4848.5713
$ws.Range("K126").Value = 31247.25
$ws.Range("L126").Value = 14545.7139
$ws.Range("M126").Value = -28777.25
$ws.Range("N126").Value = -19485.7139
$ws.Range("H132").Value = 32964.324
$ws.Range("I132").Value = 40014.78
$ws.Range("K132").Value = 120044.34
$ws.Range("M132").Value = -117514.34
$ws.Range("H136").Value = 3974
$ws.Range("I136").Value = 3974
$ws.Range("K136").Value = 11922
$ws.Range("M136").Value = -9372

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 113277.5
$ws.Range("I126").Value = 202308.5
$ws.Range("J126").Value = 1988.75
$ws.Range("K126").Value = 606925.5
$ws.Range("L126").Value = 5966.25
$ws.Range("M126").Value = -604455.5
$ws.Range("N126").Value = -10906.25
$ws.Range("H136").Value = 4250.25
$ws.Range("I136").Value = 4195
$ws.Range("J136").Value = 4447.5713
$ws.Range("K136").Value = 12585
$ws.Range("L136").Value = 13342.7139
$ws.Range("M136").Value = -10035
$ws.Range("N136").Value = -18442.7139

Write-Host "Applied all market price/profit updates"